{"js": "// \"Complemento Mec\u00e2nica da batalha\"\n// 1) Extend the \"Ao come\u00e7ar a partida...\" paragraph with two more sentences\n//    about effect cards and trap cards.\n// 2) Add a brand-new closing paragraph explaining the win/lose condition.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that starts the \"Batalha\" card-mechanics description\n// (the one beginning with \"Ao come\u00e7ar a partida\") so the script is robust\n// even if paragraph indices shift.\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Ao come\u00e7ar a partida\") === 0) {\n    targetPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetPara) {\n  throw new Error(\"Could not find the 'Ao come\u00e7ar a partida' paragraph.\");\n}\n\n// 1) Append the two new sentences to the end of that paragraph.\nconst addition =\n  \" As cartas de efeito servem para causar algum efeito em algum componente do jogo (por exemplo, aumentar o dano que os fighters causam no inimigo). \" +\n  \"J\u00e1 as cartas armadilha ficam esperando que algum jogador inimigo passe por cima dela, causando assim algum efeito no mesmo (como reduzir a velocidade da movimenta\u00e7\u00e3o ou reduzir a quantidade de vida por exemplo).\";\n\ntargetPara.insertText(addition, Word.InsertLocation.end);\n\n// 2) Insert a brand-new paragraph right after it with the win/lose wrap-up.\nconst newParaText =\n  \"O Jogador dever\u00e1 utilizar as cartas, criando uma estrat\u00e9gia para destruir a base inimiga sem deixar que o oponente destrua sua base. \" +\n  \"Uma partida normal ter\u00e1 fim quando uma das bases for destru\u00edda. Em algumas batalhas, condi\u00e7\u00f5es especiais encerrar\u00e3o a batalha. \" +\n  \"Se a base destru\u00edda for a base inimiga, o jogador ganha o jogo. Caso contr\u00e1rio, ele perde.\";\n\ntargetPara.insertParagraph(newParaText, Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# \"Complemento Mec\u00e2nica da batalha\"\n# 1) Extend the \"Ao come\u00e7ar a partida...\" paragraph with two more sentences\n#    about effect cards and trap cards.\n# 2) Add a brand-new closing paragraph explaining the win/lose condition.\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that starts the \"Batalha\" card-mechanics description\n# (the one beginning with \"Ao come\u00e7ar a partida\") so the script is robust\n# even if paragraph indices shift. Compare against an ASCII-only prefix to\n# sidestep any accent-encoding surprises in the literal string match.\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Ao come\")) {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -eq $null) {\n    throw \"Could not find the 'Ao come\u00e7ar a partida' paragraph.\"\n}\n\n$targetIndex = $targetPara.Index\n\n# 1) Append the two new sentences to the end of that paragraph.\n$addition = \" As cartas de efeito servem para causar algum efeito em algum componente do jogo (por exemplo, aumentar o dano que os fighters causam no inimigo). J\u00e1 as cartas armadilha ficam esperando que algum jogador inimigo passe por cima dela, causando assim algum efeito no mesmo (como reduzir a velocidade da movimenta\u00e7\u00e3o ou reduzir a quantidade de vida por exemplo).\"\n$targetPara.Range.InsertAfter($addition)\n\n# 2) Insert a brand-new paragraph right after it with the win/lose wrap-up.\n$newParaText = \"O Jogador dever\u00e1 utilizar as cartas, criando uma estrat\u00e9gia para destruir a base inimiga sem deixar que o oponente destrua sua base. Uma partida normal ter\u00e1 fim quando uma das bases for destru\u00edda. Em algumas batalhas, condi\u00e7\u00f5es especiais encerrar\u00e3o a batalha. Se a base destru\u00edda for a base inimiga, o jogador ganha o jogo. Caso contr\u00e1rio, ele perde.\"\n$targetPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetIndex + 1)\n$newPara.Range.Text = $newParaText\n"}
